$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Materias primas" (raw materials) column: each ingredient now
# carries a numeric quantity prefix (e.g. "2-harina" instead of "harina").
$ws.Range("C2").Value = "2-harina,1-huevos,1-vainilla,5-leche"
$ws.Range("C3").Value = "1-harina,2-manzana,5-huevos"
$ws.Range("C4").Value = "2-huevos,5-harina,1-vainilla"
$ws.Range("C5").Value = "5-harina,5-huevos"
$ws.Range("C6").Value = "1-crema,2-limon,5-merengue,4-harina,5-huevos"

# Leave the selection where the author ended up after editing the sheet.
[void]$ws.Range("D9").Select()
